$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" "54.586.90"
Set-TextCell $ws "E2" "  +5.35%  "

# Row 3
Set-TextCell $ws "D3" "3.177.88"
Set-TextCell $ws "E3" "  +2.21%  "

# Row 4
Set-TextCell $ws "E4" "  +0.04%  "

# Row 5
Set-TextCell $ws "D5" "401.42"
Set-TextCell $ws "E5" "  +3.29%  "

# Row 6
Set-TextCell $ws "D6" "109.43"
Set-TextCell $ws "E6" "  +5.48%  "

# Row 7
Set-TextCell $ws "D7" "0.551"
Set-TextCell $ws "E7" "  +1.04%  "

# Row 8
Set-TextCell $ws "D8" "0.999"
Set-TextCell $ws "E8" "  -0.05%  "

# Row 9
Set-TextCell $ws "E9" "  +4.43%  "

# Row 10
Set-TextCell $ws "D10" "39.17"
Set-TextCell $ws "E10" "  +5.26%  "

# Row 11
Set-TextCell $ws "D11" "0.0888"
Set-TextCell $ws "E11" "  +3.06%  "

# Row 12
Set-TextCell $ws "E12" "  +1.71%  "

# Row 13
Set-TextCell $ws "D13" "3.686.75"
Set-TextCell $ws "E13" "  +2.45%  "

# Row 14
Set-TextCell $ws "B14" "Polkadot"
Set-TextCell $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D14" "8.09"
Set-TextCell $ws "E14" "  +2.34%  "

# Row 15
Set-TextCell $ws "B15" "Chainlink"
Set-TextCell $ws "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D15" "19.07"
Set-TextCell $ws "E15" "  +1.75%  "

# Row 16
Set-TextCell $ws "E16" "  +7.84%  "

# Row 17
Set-TextCell $ws "D17" "3.178.03"
Set-TextCell $ws "E17" "  +2.15%  "

# Row 18
Set-TextCell $ws "D18" "10.54"
Set-TextCell $ws "E18" "  -1.48%  "

# Row 19
Set-TextCell $ws "D19" "54.492.56"
Set-TextCell $ws "E19" "  +4.98%  "

# Row 20
Set-TextCell $ws "D20" "3.29"
Set-TextCell $ws "E20" "  +2.77%  "

# Row 21
Set-TextCell $ws "D21" "0.0000101"
Set-TextCell $ws "E21" "  +3.96%  "

# Row 22
Set-TextCell $ws "D22" "12.88"
Set-TextCell $ws "E22" "  +2.95%  "

# Row 23
Set-TextCell $ws "D23" "72.34"
Set-TextCell $ws "E23" "  +3.33%  "

# Row 24
Set-TextCell $ws "D24" "275.40"
Set-TextCell $ws "E24" "  +2.49%  "

# Row 25
Set-TextCell $ws "D25" "3.26"
Set-TextCell $ws "E25" "  +4.30%  "

# Row 26
Set-TextCell $ws "D26" "7.98"
Set-TextCell $ws "E26" "  -1.56%  "

# Row 27
Set-TextCell $ws "B27" "EthereumClassic"
Set-TextCell $ws "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D27" "27.89"
Set-TextCell $ws "E27" "  +2.72%  "

# Row 28
Set-TextCell $ws "B28" "RenderToken"
Set-TextCell $ws "C28" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D28" "7.60"
Set-TextCell $ws "E28" "  +5.42%  "

# Row 29
Set-TextCell $ws "B29" "Dai"
Set-TextCell $ws "C29" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D29" "0.999"
Set-TextCell $ws "E29" "  -0.08%  "

# Row 30
Set-TextCell $ws "B30" "Kaspa"
Set-TextCell $ws "C30" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws "D30" "0.168"
Set-TextCell $ws "E30" "  -1.11%  "

# Row 31
Set-TextCell $ws "D31" "0.112"
Set-TextCell $ws "E31" "  +2.64%  "

# Row 32
Set-TextCell $ws "D32" "11.03"
Set-TextCell $ws "E32" "  +6.28%  "

# Row 33
Set-TextCell $ws "E33" "  +12.76%  "

# Row 34
Set-TextCell $ws "D34" "36.82"
Set-TextCell $ws "E34" "  +3.72%  "

# Row 35
Set-TextCell $ws "E35" "  +1.43%  "

# Row 36
Set-TextCell $ws "D36" "51.36"
Set-TextCell $ws "E36" "  +1.84%  "

# Row 37
Set-TextCell $ws "D37" "3.63"
Set-TextCell $ws "E37" "  +6.20%  "

# Row 38
Set-TextCell $ws "E38" "  -0.13%  "

# Row 39
Set-TextCell $ws "D39" "2.88"
Set-TextCell $ws "E39" "  +10.81%  "

# Row 40
Set-TextCell $ws "D40" "4.09"
Set-TextCell $ws "E40" "  +10.63%  "

# Row 41
Set-TextCell $ws "E41" "  +2.60%  "

# Row 42
Set-TextCell $ws "D42" "0.291"
Set-TextCell $ws "E42" "  +1.01%  "

# Row 43
Set-TextCell $ws "D43" "17.29"
Set-TextCell $ws "E43" "  +2.37%  "

# Row 44
Set-TextCell $ws "D44" "131.24"
Set-TextCell $ws "E44" "  +1.83%  "

# Row 45
Set-TextCell $ws "E45" "  +1.06%  "

# Row 46
Set-TextCell $ws "D46" "22.09"
Set-TextCell $ws "E46" "  -1.11%  "

# Row 47
Set-TextCell $ws "E47" "  -2.60%  "

# Row 48
Set-TextCell $ws "E48" "  -0.57%  "

# Row 49
Set-TextCell $ws "D49" "2.092.02"
Set-TextCell $ws "E49" "  +2.10%  "

# Row 50
Set-TextCell $ws "D50" "0.0346"
Set-TextCell $ws "E50" "  +8.71%  "

# Row 51
Set-TextCell $ws "D51" "0.0505"
Set-TextCell $ws "E51" "  +10.60%  "
